# PrecioFrutaHortalizas - Hortaliza, Agrícola del Norte S.A. de Arica - Coliflor
# "Fruta / hortaliza, semanal" weekly refresh:
#   Two new weekly price rows (Primera / Segunda, fecha 44523) are inserted
#   right before the former row 62, pushing the rest of the table (old rows
#   62-79) down by two rows (to 64-81). The sheet's used range grows from
#   A1:R79 to A1:R81.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the old row 62 (shifts old rows 62..79 -> 64..81)
$ws.Range("A62:A63").EntireRow.Insert()

# --- New row 62: Coliflor, Primera ---
$ws.Range("A62").Value = 1
$ws.Range("B62").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C62").Value = "Arica y Parinacota"
$ws.Range("D62").Value = 44523
$ws.Range("E62").Value = 15
$ws.Range("F62").Value = 100112008
$ws.Range("G62").Value = "Coliflor"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 1000
$ws.Range("K62").Value = 600
$ws.Range("L62").Value = 700
$ws.Range("M62").Value = 650
$ws.Range("N62").Value = "$/unidad"
$ws.Range("O62").Value = "Región de Arica y Parinacota"
$ws.Range("P62").Value = 650
$ws.Range("Q62").Value = 1
$ws.Range("R62").Value = "Hortaliza"

# --- New row 63: Coliflor, Segunda ---
$ws.Range("A63").Value = 1
$ws.Range("B63").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C63").Value = "Arica y Parinacota"
$ws.Range("D63").Value = 44523
$ws.Range("E63").Value = 15
$ws.Range("F63").Value = 100112008
$ws.Range("G63").Value = "Coliflor"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Segunda"
$ws.Range("J63").Value = 1200
$ws.Range("K63").Value = 400
$ws.Range("L63").Value = 500
$ws.Range("M63").Value = 450
$ws.Range("N63").Value = "$/unidad"
$ws.Range("O63").Value = "Región de Arica y Parinacota"
$ws.Range("P63").Value = 450
$ws.Range("Q63").Value = 1
$ws.Range("R63").Value = "Hortaliza"
